$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5254.6665
$ws.Range("I62").Value = 1721
$ws.Range("J62").Value = 6264.2856
$ws.Range("K62").Value = 1721
$ws.Range("L62").Value = 6264.2856
$ws.Range("M62").Value = -1097
$ws.Range("N62").Value = -7512.2856

$ws.Range("H64").Value = 4882.353

$ws.Range("H65").Value = 5254.6665
$ws.Range("I65").Value = 1721
$ws.Range("J65").Value = 6264.2856
$ws.Range("K65").Value = 8605
$ws.Range("L65").Value = 31321.428
$ws.Range("M65").Value = -5485
$ws.Range("N65").Value = -37561.428

$ws.Range("H67").Value = 4882.353

$ws.Range("H69").Value = 18911
$ws.Range("I69").Value = 12999.667
$ws.Range("J69").Value = 22457.8
$ws.Range("K69").Value = 38999.001
$ws.Range("L69").Value = 67373.39999999999
$ws.Range("M69").Value = -38125.001
$ws.Range("N69").Value = -69121.39999999999

$ws.Range("H70").Value = 3208.4
$ws.Range("I70").Value = 1647
$ws.Range("J70").Value = 3598.75
$ws.Range("K70").Value = 4941
$ws.Range("L70").Value = 10796.25
$ws.Range("M70").Value = -4671
$ws.Range("N70").Value = -11336.25

$ws.Range("H72").Value = 18911
$ws.Range("I72").Value = 12999.667
$ws.Range("J72").Value = 22457.8
$ws.Range("K72").Value = 116997.003
$ws.Range("L72").Value = 202120.2
$ws.Range("M72").Value = -112629.003
$ws.Range("N72").Value = -210856.2

$ws.Range("H73").Value = 3208.4
$ws.Range("I73").Value = 1647
$ws.Range("J73").Value = 3598.75
$ws.Range("K73").Value = 4941
$ws.Range("L73").Value = 10796.25
$ws.Range("M73").Value = -4005
$ws.Range("N73").Value = -12668.25

$ws.Range("H100").Value = 3578.7778
$ws.Range("J100").Value = 3752.1428
$ws.Range("L100").Value = 3752.1428
$ws.Range("N100").Value = -4834.1428

$ws.Range("H137").Value = 6175.9287
$ws.Range("I137").Value = 3423.25
$ws.Range("J137").Value = 9846.166999999999
$ws.Range("K137").Value = 10269.75
$ws.Range("L137").Value = 29538.501
$ws.Range("M137").Value = -7719.75
$ws.Range("N137").Value = -34638.501

$ws.Range("H138").Value = 1451374.1
$ws.Range("I138").Value = 920.5454999999999
$ws.Range("J138").Value = 2780956.5
$ws.Range("K138").Value = 2761.6365
$ws.Range("L138").Value = 8342869.5
$ws.Range("M138").Value = 2378.3635
$ws.Range("N138").Value = -8353149.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2044.1052
$ws.Range("I45").Value = 973.5
$ws.Range("J45").Value = 2538.2307
$ws.Range("K45").Value = 973.5
$ws.Range("L45").Value = 2538.2307
$ws.Range("M45").Value = -596.5
$ws.Range("N45").Value = -3292.2307

$ws.Range("H61").Value = 45551590
$ws.Range("I61").Value = 125000850
$ws.Range("J61").Value = 152013.72
$ws.Range("K61").Value = 125000850
$ws.Range("L61").Value = 152013.72
$ws.Range("M61").Value = -125000638
$ws.Range("N61").Value = -152437.72

$ws.Range("H74").Value = 6416949
$ws.Range("I74").Value = 9260899
$ws.Range("J74").Value = 18062.916
$ws.Range("K74").Value = 9260899
$ws.Range("L74").Value = 18062.916
$ws.Range("M74").Value = -9260025
$ws.Range("N74").Value = -19810.916

$ws.Range("H77").Value = 6416949
$ws.Range("I77").Value = 9260899
$ws.Range("J77").Value = 18062.916
$ws.Range("K77").Value = 46304495
$ws.Range("L77").Value = 90314.58
$ws.Range("M77").Value = -46300127
$ws.Range("N77").Value = -99050.58

$ws.Range("H88").Value = 1713.2142
$ws.Range("I88").Value = 1460.75
$ws.Range("J88").Value = 2049.8333
$ws.Range("K88").Value = 1460.75
$ws.Range("L88").Value = 2049.8333
$ws.Range("M88").Value = -1054.75
$ws.Range("N88").Value = -2861.8333

$ws.Range("H91").Value = 1713.2142
$ws.Range("I91").Value = 1460.75
$ws.Range("J91").Value = 2049.8333
$ws.Range("K91").Value = 1460.75
$ws.Range("L91").Value = 2049.8333
$ws.Range("M91").Value = -56.75
$ws.Range("N91").Value = -4857.8333

$ws.Range("H97").Value = 1843.3334
$ws.Range("I97").Value = 1612
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1612
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1116
$ws.Range("N97").Value = -3992

$ws.Range("H102").Value = 11146.4
$ws.Range("I102").Value = 13745.272
$ws.Range("K102").Value = 13745.272
$ws.Range("M102").Value = -12123.272

$ws.Range("H132").Value = 7755.095
$ws.Range("I132").Value = 4387
$ws.Range("K132").Value = 13161
$ws.Range("M132").Value = -10631

$ws.Range("H136").Value = 45551590
$ws.Range("I136").Value = 125000850
$ws.Range("J136").Value = 152013.72
$ws.Range("K136").Value = 375002550
$ws.Range("L136").Value = 456041.16
$ws.Range("M136").Value = -375000000
$ws.Range("N136").Value = -461141.16

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1928.3334
$ws.Range("J80").Value = 1928.3334
$ws.Range("L80").Value = 1928.3334
$ws.Range("N80").Value = -3924.3334

$ws.Range("H81").Value = 42148.6
$ws.Range("J81").Value = 42148.6
$ws.Range("L81").Value = 42148.6
$ws.Range("N81").Value = -44270.6

$ws.Range("H83").Value = 1928.3334
$ws.Range("J83").Value = 1928.3334
$ws.Range("L83").Value = 9641.666999999999
$ws.Range("N83").Value = -19625.667

$ws.Range("H84").Value = 42148.6
$ws.Range("J84").Value = 42148.6
$ws.Range("L84").Value = 126445.8
$ws.Range("N84").Value = -137053.8

$ws.Range("H94").Value = 2026.5883
$ws.Range("I94").Value = 1788.0834
$ws.Range("K94").Value = 1788.0834
$ws.Range("M94").Value = -1337.0834

$ws.Range("H138").Value = 50375
$ws.Range("J138").Value = 50375
$ws.Range("L138").Value = 50375
$ws.Range("N138").Value = -60655

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 4500
$ws.Range("K4").Value = 4500
$ws.Range("M4").Value = -4388

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H31").Value = 899306.6
$ws.Range("I31").Value = 1523.8462
$ws.Range("J31").Value = 1797089.5
$ws.Range("K31").Value = 1523.8462
$ws.Range("L31").Value = 1797089.5
$ws.Range("M31").Value = -1228.8462
$ws.Range("N31").Value = -1797679.5

$ws.Range("H34").Value = 899306.6
$ws.Range("I34").Value = 1523.8462
$ws.Range("J34").Value = 1797089.5
$ws.Range("K34").Value = 1523.8462
$ws.Range("L34").Value = 1797089.5
$ws.Range("M34").Value = -1321.8462
$ws.Range("N34").Value = -1797493.5

$ws.Range("H56").Value = 17024
$ws.Range("I56").Value = 1093
$ws.Range("J56").Value = 24989.5
$ws.Range("K56").Value = 1093
$ws.Range("L56").Value = 24989.5
$ws.Range("M56").Value = -248
$ws.Range("N56").Value = -26679.5

$ws.Range("H86").Value = 2866.4285
$ws.Range("I86").Value = 3124.75
$ws.Range("K86").Value = 3124.75
$ws.Range("M86").Value = -2001.75

$ws.Range("H89").Value = 2866.4285
$ws.Range("I89").Value = 3124.75
$ws.Range("K89").Value = 15623.75
$ws.Range("M89").Value = -10007.75

$ws.Range("H107").Value = 772.7646999999999
$ws.Range("I107").Value = 501.18182
$ws.Range("J107").Value = 1270.6666
$ws.Range("K107").Value = 501.18182
$ws.Range("L107").Value = 1270.6666
$ws.Range("M107").Value = 1418.81818
$ws.Range("N107").Value = -5110.6666

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H118").Value = 87969
$ws.Range("J118").Value = 87969
$ws.Range("L118").Value = 87969
$ws.Range("N118").Value = -91283

$ws.Range("H132").Value = 2109.5334
$ws.Range("I132").Value = 1903.0714
$ws.Range("K132").Value = 5709.2142
$ws.Range("M132").Value = -3179.2142

$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 31096.23
$ws.Range("I2").Value = 395.1111
$ws.Range("J2").Value = 100173.75
$ws.Range("K2").Value = 2370.6666
$ws.Range("L2").Value = 601042.5
$ws.Range("M2").Value = -2257.6666
$ws.Range("N2").Value = -601268.5

$ws.Range("H11").Value = 324.53845
$ws.Range("I11").Value = 251.58333
$ws.Range("K11").Value = 754.74999
$ws.Range("M11").Value = -614.74999

$ws.Range("H34").Value = 3302.9412
$ws.Range("J34").Value = 4863.727
$ws.Range("L34").Value = 14591.181
$ws.Range("N34").Value = -14759.181

$ws.Range("H55").Value = 9312.4375
$ws.Range("I55").Value = 9333
$ws.Range("J55").Value = 9307.691999999999
$ws.Range("K55").Value = 27999
$ws.Range("L55").Value = 27923.076
$ws.Range("M55").Value = -27822
$ws.Range("N55").Value = -28277.076

$ws.Range("H131").Value = 5578.067
$ws.Range("J131").Value = 21766
$ws.Range("L131").Value = 65298
$ws.Range("N131").Value = -75378

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 75500
$ws.Range("J15").Value = 75500
$ws.Range("L15").Value = 75500
$ws.Range("N15").Value = -76076

$ws.Range("H81").Value = 75500
$ws.Range("J81").Value = 75500
$ws.Range("L81").Value = 75500
$ws.Range("N81").Value = -77496

$ws.Range("H84").Value = 75500
$ws.Range("J84").Value = 75500
$ws.Range("L84").Value = 226500
$ws.Range("N84").Value = -236484

$ws.Range("H122").Value = 3831
$ws.Range("I122").Value = 3197.2
$ws.Range("K122").Value = 9591.599999999999
$ws.Range("M122").Value = -7141.599999999999

$ws.Range("H132").Value = 142861550
$ws.Range("I132").Value = 200004580
$ws.Range("K132").Value = 600013740
$ws.Range("M132").Value = -600011210

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1457.3125
$ws.Range("I16").Value = 1652
$ws.Range("J16").Value = 613.6667
$ws.Range("K16").Value = 1652
$ws.Range("L16").Value = 613.6667
$ws.Range("M16").Value = -1482
$ws.Range("N16").Value = -953.6667

$ws.Range("H22").Value = 2184.55
$ws.Range("I22").Value = 2230.6875
$ws.Range("K22").Value = 2230.6875
$ws.Range("M22").Value = -1935.6875

$ws.Range("H27").Value = 2184.55
$ws.Range("I27").Value = 2230.6875
$ws.Range("K27").Value = 2230.6875
$ws.Range("M27").Value = -2123.6875

$ws.Range("H122").Value = 5190.2256
$ws.Range("I122").Value = 4691.2607
$ws.Range("K122").Value = 14073.7821
$ws.Range("M122").Value = -11623.7821

$ws.Range("H132").Value = 382353.2
$ws.Range("I132").Value = 771961.4399999999
$ws.Range("K132").Value = 2315884.32
$ws.Range("M132").Value = -2313354.32

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 157999.6
$ws.Range("I2").Value = 157999.6
$ws.Range("K2").Value = 157999.6
$ws.Range("M2").Value = -157887.6

$ws.Range("H68").Value = 41800
$ws.Range("J68").Value = 41800
$ws.Range("L68").Value = 41800
$ws.Range("N68").Value = -43422

$ws.Range("H71").Value = 41800
$ws.Range("J71").Value = 41800
$ws.Range("L71").Value = 125400
$ws.Range("N71").Value = -133512

$ws.Range("H81").Value = 41040
$ws.Range("I81").Value = 1250
$ws.Range("K81").Value = 2500
$ws.Range("M81").Value = -1439

$ws.Range("H84").Value = 41040
$ws.Range("I84").Value = 1250
$ws.Range("K84").Value = 12500
$ws.Range("M84").Value = -7196

$ws.Range("H96").Value = 1348
$ws.Range("I96").Value = 1348
$ws.Range("K96").Value = 1348
$ws.Range("M96").Value = 25

$ws.Range("H116").Value = 97990
$ws.Range("J116").Value = 97990
$ws.Range("L116").Value = 97990
$ws.Range("N116").Value = -107168

$ws.Range("H122").Value = 5241.6665
$ws.Range("I122").Value = 3708.625
$ws.Range("J122").Value = 7471.5454
$ws.Range("K122").Value = 11125.875
$ws.Range("L122").Value = 22414.6362
$ws.Range("M122").Value = -8675.875
$ws.Range("N122").Value = -27314.6362

$ws.Range("H132").Value = 3974.3428
$ws.Range("I132").Value = 3572.84
$ws.Range("J132").Value = 4978.1
$ws.Range("K132").Value = 10718.52
$ws.Range("L132").Value = 14934.3
$ws.Range("M132").Value = -8188.52
$ws.Range("N132").Value = -19994.3

$ws.Range("H136").Value = 1817.5454
$ws.Range("I136").Value = 1713.2858
$ws.Range("K136").Value = 5139.857400000001
$ws.Range("M136").Value = -2589.857400000001
